# Weekly fruit/vegetable price data update.
# The data rows (2-14) were reordered; this script rewrites the
# per-row values for columns D (Fecha), J (Volumen), K (Precio minimo),
# L (Precio maximo), M (Precio promedio ponderado), O (Origen) and
# P (Precio $/Kg) to match the new row order. Other columns
# (A,B,C,E,F,G,H,I,N,Q,R) are identical across all rows and stay untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ Row=2;  D=44446; J=25; K=14000; L=14000; M=14000; O="Provincia de Limarí"; P=467 },
    @{ Row=3;  D=44421; J=25; K=15000; L=16000; M=15400; O="Provincia de Limarí"; P=513 },
    @{ Row=4;  D=44435; J=25; K=14000; L=14000; M=14000; O="Provincia de Limarí"; P=467 },
    @{ Row=5;  D=44435; J=25; K=14000; L=14000; M=14000; O="Provincia del Elquí"; P=467 },
    @{ Row=6;  D=44460; J=45; K=13000; L=13000; M=13000; O="Provincia de Limarí"; P=433 },
    @{ Row=7;  D=44418; J=30; K=15000; L=15000; M=15000; O="Provincia de Limarí"; P=500 },
    @{ Row=8;  D=44474; J=45; K=10000; L=10000; M=10000; O="Provincia de Limarí"; P=333 },
    @{ Row=9;  D=44425; J=35; K=14000; L=14000; M=14000; O="Provincia de Limarí"; P=467 },
    @{ Row=10; D=44453; J=50; K=12000; L=12000; M=12000; O="Provincia de Limarí"; P=400 },
    @{ Row=11; D=44467; J=35; K=12000; L=12000; M=12000; O="Provincia de Limarí"; P=400 },
    @{ Row=12; D=44449; J=45; K=12000; L=12000; M=12000; O="Provincia de Limarí"; P=400 },
    @{ Row=13; D=44376; J=25; K=18000; L=18000; M=18000; O="Provincia de Limarí"; P=600 },
    @{ Row=14; D=44432; J=25; K=14000; L=14000; M=14000; O="Provincia del Elquí"; P=467 }
)

foreach ($r in $rows) {
    $row = $r.Row
    $ws.Cells.Item($row, 4).Value  = $r.D   # D: Fecha
    $ws.Cells.Item($row, 10).Value = $r.J   # J: Volumen
    $ws.Cells.Item($row, 11).Value = $r.K   # K: Precio minimo
    $ws.Cells.Item($row, 12).Value = $r.L   # L: Precio maximo
    $ws.Cells.Item($row, 13).Value = $r.M   # M: Precio promedio ponderado
    $ws.Cells.Item($row, 15).Value = $r.O   # O: Origen
    $ws.Cells.Item($row, 16).Value = $r.P   # P: Precio $/Kg
}
